$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) data range to text format so numeric-looking
# strings (e.g. "229.09", "58.90") keep their exact original text,
# instead of being auto-converted to numbers by Excel.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "41.644.22"
$ws.Range("E2").Value = "  +5.32%  "
$ws.Range("D3").Value = "2.233.17"
$ws.Range("E3").Value = "  +3.34%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "229.09"
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("D6").Value = "0.624"
$ws.Range("E6").Value = "  -1.88%  "
$ws.Range("D7").Value = "61.89"
$ws.Range("E7").Value = "  -2.67%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").Value = "0.402"
$ws.Range("E9").Value = "  +2.15%  "
$ws.Range("D10").Value = "58.90"
$ws.Range("E10").Value = "  +1.38%  "
$ws.Range("D11").Value = "0.0879"
$ws.Range("E11").Value = "  +3.25%  "
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("D13").Value = "2.563.69"
$ws.Range("E13").Value = "  +3.26%  "
$ws.Range("D14").Value = "15.69"
$ws.Range("E14").Value = "  -1.91%  "
$ws.Range("D15").Value = "22.08"
$ws.Range("E15").Value = "  +0.28%  "
$ws.Range("E16").Value = "  -0.89%  "
$ws.Range("D17").Value = "5.58"
$ws.Range("E17").Value = "  +1.42%  "
$ws.Range("D18").Value = "2.232.37"
$ws.Range("E18").Value = "  +3.33%  "
$ws.Range("D19").Value = "41.511.66"
$ws.Range("E19").Value = "  +4.75%  "
$ws.Range("D20").Value = "73.24"
$ws.Range("E20").Value = "  +1.30%  "
$ws.Range("D21").Value = "0.0₃0904"
$ws.Range("E21").Value = "  +6.75%  "
$ws.Range("E22").Value = "  -2.34%  "
$ws.Range("D23").Value = "247.62"
$ws.Range("E23").Value = "  +7.75%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "2.39"
$ws.Range("E25").Value = "  -3.95%  "
$ws.Range("D26").Value = "2.36"
$ws.Range("E26").Value = "  +0.66%  "
$ws.Range("D27").Value = "9.57"
$ws.Range("E27").Value = "  +0.62%  "
$ws.Range("D28").Value = "0.143"
$ws.Range("E28").Value = "  +2.99%  "
$ws.Range("D29").Value = "168.49"
$ws.Range("E29").Value = "  -2.15%  "
$ws.Range("D30").Value = "20.08"
$ws.Range("E30").Value = "  +1.16%  "
$ws.Range("E31").Value = "  +0.54%  "
$ws.Range("D32").Value = "2.81"
$ws.Range("E32").Value = "  +5.52%  "
$ws.Range("D33").Value = "0.122"
$ws.Range("E33").Value = "  -0.54%  "
$ws.Range("D34").Value = "4.98"
$ws.Range("E34").Value = "  +6.13%  "
$ws.Range("D35").Value = "4.62"
$ws.Range("E35").Value = "  +0.72%  "
$ws.Range("D36").Value = "0.0626"
$ws.Range("E36").Value = "  +1.27%  "
$ws.Range("E37").Value = "  +3.73%  "
$ws.Range("D38").Value = "6.69"
$ws.Range("E38").Value = "  -4.65%  "
$ws.Range("D39").Value = "2.38"
$ws.Range("E39").Value = "  -1.96%  "
$ws.Range("E40").Value = "  +0.18%  "
$ws.Range("D41").Value = "0.000239"
$ws.Range("E41").Value = "  +27.54%  "
$ws.Range("D42").Value = "4.92"
$ws.Range("E42").Value = "  +6.81%  "
$ws.Range("D43").Value = "0.0236"
$ws.Range("E43").Value = "  +4.05%  "
$ws.Range("D44").Value = "8.67"
$ws.Range("E44").Value = "  +12.43%  "
$ws.Range("D45").Value = "99.85"
$ws.Range("E45").Value = "  -2.57%  "
$ws.Range("D46").Value = "0.0958"
$ws.Range("E46").Value = "  +3.65%  "
$ws.Range("D47").Value = "1.485.54"
$ws.Range("E47").Value = "  -2.68%  "
$ws.Range("E48").Value = "  -1.24%  "
$ws.Range("D49").Value = "16.79"
$ws.Range("E49").Value = "  -5.98%  "
$ws.Range("D50").Value = "2.78"
$ws.Range("E50").Value = "  -1.03%  "
$ws.Range("E51").Value = "  -2.09%  "

# Restore the default cell style on the price range so no stray
# number-format style is left attached to these cells.
$priceRange.Style = "Normal"

